$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert new row 35 data (Literature Review entry dated 17.10.2022)
$ws.Range("A35").Value = "17.10.2022"
$ws.Range("B35").Value = 0.91666666666666663
$ws.Range("C35").Value = "Literature Review"
$ws.Range("D35").Value = "Documentation"
$ws.Range("E35").Value = 180
$ws.Range("F35").Value = "Music Notations, Terminology, Keyboard Ghosting and Diagrams"

# Copy formatting from the row above (A34/B34) so A35/B35 pick up the same
# date/time cell styles (gray fill + center alignment, time number format)
$ws.Range("A34").Copy()
$ws.Range("A35").PasteSpecial(-4122)

$ws.Range("B34").Copy()
$ws.Range("B35").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Extend totals formulas to include the new row
$ws.Range("E39").Formula = "=SUM(E2:E35)"
$ws.Range("E40").Formula = "=E39 / 60"

# Update the active selection to match the saved workbook state
$ws.Range("F37").Select()
